$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- New row 12: geopoint / location / Record your location ---
$ws.Range("A12").Value = "geopoint"
$ws.Range("B12").Value = "location"
$ws.Range("C12").Value = "Record your location"

# --- New row 13: rating_button / agreement / question ---
$ws.Range("A13").Value = "rating_button"
$ws.Range("B13").Value = "agreement"
$ws.Range("C13").Value = "Do you agree that SurveySignal is cool & hip?"

# --- Header row additions (D1, E1) ---
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "choice1"
$ws.Range("E1").Value = "choice2"

# --- Row 13 icons ---
$ws.Range("D13").Value = '<i class="fa fa-smile-o fa-2x"></i>'
$ws.Range("E13").Value = '<i class="fa fa-meh-o fa-2x"></i>'

$ws.Rows.Item(13).RowHeight = 45

# --- Selection update ---
$ws.Range("E13").Select()
